$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-23 Tuesday", "2024-04-24 Wednesday"),
    @("374÷2=", "225÷2="),
    @("161÷3=", "605÷8="),
    @("213÷6=", "678÷2="),
    @("846÷9=", "761÷9="),
    @("297÷5=", "734÷6="),
    @("337÷5=", "667÷9="),
    @("664÷4=", "858÷2="),
    @("756÷9=", "136÷5="),
    @("637÷4=", "600÷2="),
    @("407÷2=", "177÷9="),
    @("722÷3=", "356÷5="),
    @("390÷7=", "645÷8="),
    @("227÷6=", "524÷9="),
    @("198÷6=", "191÷8="),
    @("242÷9=", "939÷8="),
    @("266÷9=", "543÷6="),
    @("134÷9=", "455÷2="),
    @("194÷9=", "172÷3="),
    @("793÷7=", "711÷5="),
    @("730÷2=", "560÷2="),
    @("255÷2=", "580÷7="),
    @("826÷9=", "718÷7="),
    @("897÷9=", "601÷8="),
    @("529÷3=", "835÷9="),
    @("173÷4=", "470÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
